$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log - Part 3")
$ws.Activate()

# --- Row 8 ---
$ws.Range("B8").Value = 6977
$ws.Range("C8").Value = 43930
$ws.Range("D8").Value = 0.9
$ws.Range("E8").Value = 0.91319444444444453
$ws.Range("G8").Value = "Compared constraints for Logic Unit and Arithmetic unit to our implementations. Taking a break for dinner/supper. NOT DONE"

# --- Row 9 ---
$ws.Range("B9").Value = 6977
$ws.Range("C9").Value = 43930
$ws.Range("D9").Value = 0.95833333333333337
$ws.Range("E9").Value = 0.96388888888888891
$ws.Range("G9").Value = "Checked that implementation for our LogicUnit.vhd fulfills constraints. DONE"

# --- Row 10 ---
$ws.Range("B10").Value = 6977
$ws.Range("C10").Value = 43930
$ws.Range("D10").Value = 0.96388888888888891
$ws.Range("E10").Value = 0.97152777777777777
$ws.Range("G10").Value = "Edited ArithUnit.vhd to fulfill the constraint of adding an extra output directy from the Adder. Checked that implementation for our ArithUnit.vhd fulfills constraints. DONE"

# --- Row 11 ---
$ws.Range("B11").Value = 6977
$ws.Range("C11").Value = 43930
$ws.Range("D11").Value = 0.97152777777777777
$ws.Range("E11").Value = 0.9770833333333333
$ws.Range("G11").Value = "Added in VHDL interface to SLL64.vhd, SLRA64.vhd and SRL64.vhd. DONE. Taking a small break"

# --- Row 12 (start time only; end time / description remain blank) ---
$ws.Range("B12").Value = 6977
$ws.Range("C12").Value = 43930
$ws.Range("D12").Value = 0.027083333333333334

# Move the live selection to match the author's last cursor position
$ws.Range("F12").Select()
